# Insert a new price record for "Feria Lagunitas de Puerto Montt - Arveja Verde"
# right after the existing row 15. This pushes the previously existing rows
# 16-71 down to 17-72 (dimension grows from A1:R71 to A1:R72) and populates
# the newly freed row 16 with a new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 16, shifting rows 16:71 down to 17:72
$ws.Rows("16:16").Insert()

# Populate the new row 16 with the new record
$ws.Range("A16").Value2 = 4
$ws.Range("B16").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C16").Value2 = "Los Lagos"
$ws.Range("D16").Value2 = 44481
$ws.Range("E16").Value2 = 10
$ws.Range("F16").Value2 = 100112022
$ws.Range("G16").Value2 = "Arveja Verde"
$ws.Range("H16").Value2 = "Sin especificar"
$ws.Range("I16").Value2 = "Primera"
$ws.Range("J16").Value2 = 80
$ws.Range("K16").Value2 = 25000
$ws.Range("L16").Value2 = 25000
$ws.Range("M16").Value2 = 25000
$ws.Range("N16").Value2 = "$/saco 25 kilos"
$ws.Range("O16").Value2 = "Región Metropolitana"
$ws.Range("P16").Value2 = 1000
$ws.Range("Q16").Value2 = 25
$ws.Range("R16").Value2 = "Hortaliza"
